$d = $word.ActiveDocument

# Helper: force a run split at [start,end) by toggling a no-op formatting
# change (Bold on then off). The runtime's XML writer only splits a run
# into pieces at locations where a distinct formatting write is applied,
# even though the value ends up identical to what was already there.
function Split-Run($start, $end) {
    $r = $d.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

# --- 1. Primary-color: OB3954  ->  "Primary-color: #0" | "B3954" ---------
$rng = $d.Content
$null = $rng.Find.Execute("Primary-color: OB3954", $true, $false, $false, $false, $false, $true, 1, $false, "Primary-color: #0B3954", 2)

$rng2 = $d.Content
$null = $rng2.Find.Execute("Primary-color: #0B3954")
$e2 = $rng2.End
$s2 = $e2 - 5
Split-Run $s2 $e2

# --- 2. Secondary-color: 087E8B -> "Secondary-color: " | "#" | "087E8B" --
$rng = $d.Content
$null = $rng.Find.Execute("Secondary-color: ")
$p = $rng.End
$ip = $d.Range($p, $p)
$ip.InsertAfter("#")
$pEnd = $p + 1
Split-Run $p $pEnd

# --- 3. Accent1-color: FF5A5F -> "Accent1-color: " | "#" | "FF5A5F" ------
$rng = $d.Content
$null = $rng.Find.Execute("Accent1-color: ")
$p = $rng.End
$ip = $d.Range($p, $p)
$ip.InsertAfter("#")
$pEnd = $p + 1
Split-Run $p $pEnd

# --- 4. Accent2-color: F2F2F2 -> "Accent2-color: " | "#" | <bm> | "F2F2F2"
$rng = $d.Content
$null = $rng.Find.Execute("Accent2-color: ")
$p = $rng.End
$ip = $d.Range($p, $p)
$ip.InsertAfter("#")
$pEnd = $p + 1
Split-Run $p $pEnd

# Move the _GoBack bookmark here, between the "#" run and the "F2F2F2" run.
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete()
} catch {
}
$bmRange = $d.Range($pEnd, $pEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
